# Remove the stale percentage-based "RPA status" notes that had been
# mistakenly stashed in column F (rows 3-5). This also drops the three
# now-unused shared strings and renumbers the rest automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3:F5").ClearContents()

# Match the author's final selection/view state (select the now-empty
# cells that were just cleared).
$ws.Range("F3:F5").Select()
